$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Value used by the other "Water" rows (1e-9, with negative sign as in row 14)
$val = -0.000000001

# Row 15: "Water" vs "water::fossil well"
$ws.Range("A15").Value = "Water"
$ws.Range("B15").Value = "water::fossil well"
$ws.Range("C15").Value = $val

# Row 16: "Water" vs "water::ground-, long-term"
$ws.Range("A16").Value = "Water"
$ws.Range("B16").Value = "water::ground-, long-term"
$ws.Range("C16").Value = $val

# Copy the cell formatting (highlight fill, style id 3) from the last existing
# "Water" row (row 14) onto the two new rows, matching the target workbook.
$ws.Range("A14:C14").Copy()
$ws.Range("A15:C16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Move the active selection the way Excel would after appending the rows.
$ws.Range("B19").Select()
